$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-07-08 Monday" "2024-07-09 Tuesday"

Replace-Text "908÷7=" "613÷7="
Replace-Text "886÷4=" "720÷9="
Replace-Text "127÷7=" "730÷2="
Replace-Text "616÷8=" "383÷7="
Replace-Text "558÷4=" "331÷3="
Replace-Text "149÷2=" "789÷8="
Replace-Text "583÷7=" "279÷4="
Replace-Text "994÷3=" "383÷2="
Replace-Text "645÷8=" "369÷3="
Replace-Text "489÷3=" "760÷8="
Replace-Text "671÷5=" "231÷2="
Replace-Text "281÷6=" "397÷6="
Replace-Text "687÷7=" "549÷7="
Replace-Text "913÷9=" "857÷4="
Replace-Text "899÷6=" "808÷6="
Replace-Text "614÷6=" "512÷7="
Replace-Text "129÷2=" "107÷8="
Replace-Text "334÷3=" "490÷6="
Replace-Text "384÷7=" "320÷8="
Replace-Text "164÷9=" "372÷6="
Replace-Text "494÷5=" "976÷5="
Replace-Text "598÷7=" "708÷4="
Replace-Text "416÷4=" "494÷3="
Replace-Text "137÷5=" "722÷7="
Replace-Text "628÷7=" "442÷9="
